$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $c = $ws.Range($addr)
    $c.Value = "'" + $val
    $c.Style = "Normal"
}

Set-TextCell 'D2' '69.545.03'
Set-TextCell 'E2' '  +1.35%  '
Set-TextCell 'D3' '2.441.19'
Set-TextCell 'E3' '  -0.08%  '
Set-TextCell 'E4' '  +0.21%  '
Set-TextCell 'D5' '564.19'
Set-TextCell 'E5' '  +0.78%  '
Set-TextCell 'D6' '165.75'
Set-TextCell 'E6' '  +1.54%  '
Set-TextCell 'E7' '  +0.08%  '
Set-TextCell 'D8' '0.510'
Set-TextCell 'E8' '  -0.16%  '
Set-TextCell 'E9' '  +10.64%  '
Set-TextCell 'E10' '  -1.48%  '
Set-TextCell 'E11' '  +1.86%  '
Set-TextCell 'D12' '4.66'
Set-TextCell 'E12' '  -2.96%  '
Set-TextCell 'D13' '0.0000181'
Set-TextCell 'E13' '  +6.10%  '
Set-TextCell 'D14' '69.440.82'
Set-TextCell 'E14' '  +1.46%  '
Set-TextCell 'D15' '2.893.06'
Set-TextCell 'E15' '  -0.44%  '
Set-TextCell 'D16' '23.92'
Set-TextCell 'E16' '  +2.83%  '
Set-TextCell 'D17' '2.444.60'
Set-TextCell 'E17' '  +0.91%  '
Set-TextCell 'D18' '10.77'
Set-TextCell 'E18' '  +3.48%  '
Set-TextCell 'D19' '341.10'
Set-TextCell 'E19' '  +1.08%  '
Set-TextCell 'D20' '7.08'
Set-TextCell 'E20' '  +2.91%  '
Set-TextCell 'D21' '3.87'
Set-TextCell 'E21' '  +1.86%  '
Set-TextCell 'D22' '2.00'
Set-TextCell 'E22' '  +6.20%  '
Set-TextCell 'E23' '  +0.03%  '
Set-TextCell 'D24' '66.09'
Set-TextCell 'E24' '  -1.12%  '
Set-TextCell 'D25' '3.86'
Set-TextCell 'E25' '  +4.57%  '
Set-TextCell 'B26' 'WrappedeETH'
Set-TextCell 'C26' 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
Set-TextCell 'D26' '2.571.54'
Set-TextCell 'E26' '  +0.34%  '
Set-TextCell 'B27' 'Aptos'
Set-TextCell 'C27' 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextCell 'D27' '8.50'
Set-TextCell 'E27' '  +4.53%  '
Set-TextCell 'D28' '0.997'
Set-TextCell 'E28' '  -0.38%  '
Set-TextCell 'D29' '0.0₃0850'
Set-TextCell 'E29' '  +4.39%  '
Set-TextCell 'D30' '7.30'
Set-TextCell 'E30' '  +1.46%  '
Set-TextCell 'D31' '1.24'
Set-TextCell 'E31' '  +8.54%  '
Set-TextCell 'D32' '447.05'
Set-TextCell 'E32' '  +4.90%  '
Set-TextCell 'D33' '1.00'
Set-TextCell 'E33' '  +0.16%  '
Set-TextCell 'E34' '  +0.40%  '
Set-TextCell 'D35' '162.12'
Set-TextCell 'E35' '  +0.70%  '
Set-TextCell 'E36' '  +0.28%  '
Set-TextCell 'E37' '  -0.02%  '
Set-TextCell 'D38' '0.109'
Set-TextCell 'E38' '  +3.39%  '
Set-TextCell 'D39' '18.09'
Set-TextCell 'E39' '  +1.66%  '
Set-TextCell 'D40' '0.305'
Set-TextCell 'E40' '  +2.97%  '
Set-TextCell 'E41' '  +4.43%  '
Set-TextCell 'D42' '4.43'
Set-TextCell 'E42' '  +1.51%  '
Set-TextCell 'E43' '  +1.81%  '
Set-TextCell 'D44' '2.13'
Set-TextCell 'E44' '  +5.80%  '
Set-TextCell 'E45' '  +0.95%  '
Set-TextCell 'D46' '131.71'
Set-TextCell 'E46' '  +1.51%  '
Set-TextCell 'E47' '  +1.10%  '
Set-TextCell 'D48' '0.488'
Set-TextCell 'E48' '  +1.61%  '
Set-TextCell 'D49' '0.560'
Set-TextCell 'E49' '  +0.08%  '
Set-TextCell 'D50' '0.0930'
Set-TextCell 'E50' '  +1.34%  '
Set-TextCell 'E51' '  +2.84%  '
